$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, E, F, G, H, L, M across rows 2-25
# (Sheet: pl_mw.xlsx - res_line results, case with 380 kV)
$newValues = @{
    2 = @(1.305921270099134, 0.2152587901791208, 0.09917698049238033, 0.4443680307746263, 1.1768764079997, 1.128060550239411, 0.191105749164862, 0.2723799715095581)
    3 = @(1.2105800307354, 0.1994966497811959, 0.09978732481234609, 0.387822817061874, 1.175406406532517, 1.134412510248467, 0.1886071104622502, 0.2579451600480311)
    4 = @(1.152527770935308, 0.1897371273770716, 0.1001838364091892, 0.3531389305168915, 1.175547655107906, 1.139043655788385, 0.1871687783696316, 0.2492056783358407)
    5 = @(1.128994022145889, 0.1857395975149814, 0.1003508996730308, 0.3390132514313251, 1.175866309931621, 1.14111419861193, 0.1866067537436464, 0.2456753893139592)
    6 = @(1.125093703362666, 0.1850745774457323, 0.1003789717804285, 0.336668177824194, 1.175934951233671, 1.141469067531432, 0.1865148862422075, 0.2450910687807593)
    7 = @(1.152209887930553, 0.1896832980009719, 0.100186067274169, 0.3529483938344953, 1.175550897278384, 1.139070838304988, 0.1871611010883925, 0.2491579415439418)
    8 = @(1.272946597811767, 0.2098409446449239, 0.09938291769764662, 0.4248636149813336, 1.176152065784933, 1.130098677394571, 0.1902243292102526, 0.2673772310109257)
    9 = @(1.513577177098398, 0.2487231174117142, 0.09798010131798585, 0.5661985755041457, 1.18567754578207, 1.118328112427434, 0.1969920609151359, 0.3040855089522907)
    10 = @(1.692742716562236, 0.2768981162372768, 0.09705372474671325, 0.6702781546542269, 1.197855598500269, 1.113264191277864, 0.2024293592131556, 0.3316561037249528)
    11 = @(1.774770753947621, 0.2896314006033265, 0.09665478927531357, 0.7176906081379002, 1.204540266743294, 1.111746017384377, 0.205004252476499, 0.3443301952085491)
    12 = @(1.805908138595839, 0.2944411438116674, 0.0965069444256591, 0.7356546913071611, 1.207237761854913, 1.111284653026019, 0.2059938970759561, 0.3491485569100092)
    13 = @(1.799198808132758, 0.293405817266148, 0.09653864224130015, 0.7317853510981394, 1.206649393711757, 1.111378957145007, 0.2057801105451063, 0.3481099938574701)
    14 = @(1.777330940935144, 0.290027343316325, 0.09664256144061811, 0.7191683204515869, 1.20475885247285, 1.111705781852351, 0.2050853787473272, 0.3447262249636793)
    15 = @(1.763946018523882, 0.2879563566741581, 0.09670663446294914, 0.7114413442032514, 1.203622525409202, 1.111920774639941, 0.2046617353664715, 0.3426560367490481)
    16 = @(1.687392523367123, 0.2760642787260394, 0.09708024766480927, 0.6671810134426437, 1.197441907287924, 1.113379262856228, 0.2022631252080203, 0.3308304768457475)
    17 = @(1.640563602251916, 0.2687474127482119, 0.09731519800382427, 0.6400460337215605, 1.193944602280197, 1.114475585448986, 0.2008176373387585, 0.3236096832307425)
    18 = @(1.6136782526915, 0.2645310724985563, 0.09745245147887749, 0.6244449056556647, 1.192040708133845, 1.115180054058754, 0.1999957787090807, 0.3194689085144589)
    19 = @(1.604583830898662, 0.2631021414604788, 0.09749928688920784, 0.6191636801734006, 1.191414526845222, 1.115431248201816, 0.1997191513032277, 0.3180690501972094)
    20 = @(1.645543509707636, 0.2695271207219605, 0.09728996817996916, 0.642933953830422, 1.194305742558697, 1.114351228156607, 0.20097052382512, 0.3243770623592113)
    21 = @(1.783752026471916, 0.2910200110618746, 0.0966119504426175, 0.7228739723492197, 1.205309628586974, 1.111606699477051, 0.2052890423736642, 0.3457196056540326)
    22 = @(1.874517172741719, 0.3049964561109846, 0.09618761074574, 0.7751780083420101, 1.213470504556483, 1.110475037205362, 0.2081964777498797, 0.3597786813909281)
    23 = @(1.8260341339452, 0.2975434196389131, 0.09641237301511685, 0.7472568307915566, 1.209025701917199, 1.111018253386931, 0.2066369433694746, 0.3522649890321645)
    24 = @(1.643291977587296, 0.269174645198774, 0.09730136780055554, 0.6416283278902313, 1.194142138631605, 1.114407219072206, 0.2009013752772546, 0.3240300973735444)
    25 = @(1.448064636123661, 0.2382734589961615, 0.09834124016102574, 0.5279251897347308, 1.182197653387277, 1.120885420712156, 0.1950796388662752, 0.2940496627697513)
}

$columns = @("B", "C", "E", "F", "G", "H", "L", "M")

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
